$d = $word.ActiveDocument

# The whole body currently has 3 paragraphs:
#   1) "Once you are inside the repository:"
#   2) "git add *    " + tab + tab + "To add all files"   (list item)
#   3) (empty, contains the _GoBack bookmark)              (list item)
#
# The target document has 7 paragraphs. We replace the full span covering
# all three original paragraphs (from the very start of paragraph 1 through
# the end of paragraph 3, i.e. the whole document body) with freshly
# authored OOXML for all seven target paragraphs in a single InsertXML
# call. Doing it in one shot (rather than touching the very last paragraph
# of the document on its own) avoids leaving a stray empty trailing
# paragraph behind.

$p1 = $d.Paragraphs(1)
$p3 = $d.Paragraphs(3)
$full = $d.Range($p1.Range.Start, $p3.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">To upload changes into the remote repository in </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>github</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t>:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Go to the directory where the repository is located</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Go inside the repository</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">git add *    </w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:t>To add all files</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">git status </w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>To</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> check the status</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>git commit -m "</w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>CaesarCipher</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> was added"</w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>To commit</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Prrafodelista"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>git pus -u origin master</w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>To</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> push into remote repository</w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/>
                <w:sz w:val="18"/>
                <w:szCs w:val="18"/>
              </w:rPr>
              <w:tab/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$full.InsertXML($xml)
